$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text so the stored value matches the literal string.
$ws.Range('D2').Value = '69.281.03'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '3.672.52'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '675.83'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('E6').Value = '  -2.31%  '
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.91'
$ws.Range('E10').Value = '  -5.56%  '
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('E12').Value = '  -3.03%  '
$ws.Range('D13').Value = '4.294.13'
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.34'
$ws.Range('E14').Value = '  -3.85%  '
$ws.Range('D15').Value = '3.666.33'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').Value = '69.226.60'
$ws.Range('E16').Value = '  -0.28%  '
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '16.04'
$ws.Range('E18').Value = '  -1.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.43'
$ws.Range('E19').Value = '  -2.77%  '
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.99'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.648'
$ws.Range('E22').Value = '  -2.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '79.75'
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('D24').Value = '3.818.66'
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -6.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.89'
$ws.Range('E27').Value = '  -5.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.06'
$ws.Range('E28').Value = '  -4.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.67'
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('E30').Value = '  -4.65%  '
$ws.Range('E31').Value = '  -3.25%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.99'
$ws.Range('E33').Value = '  -4.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.88'
$ws.Range('E34').Value = '  -0.72%  '
$ws.Range('D35').Value = '3.665.76'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('E36').Value = '  -4.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.18'
$ws.Range('E37').Value = '  -3.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.22'
$ws.Range('E38').Value = '  -1.05%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  -3.74%  '
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '173.70'
$ws.Range('E43').Value = '  +6.99%  '
$ws.Range('E44').Value = '  -1.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.57'
$ws.Range('E45').Value = '  -1.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.29'
$ws.Range('E46').Value = '  -5.64%  '
$ws.Range('B47').Value = 'FLOKI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.000278'
$ws.Range('E47').Value = '  -3.04%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.69'
$ws.Range('E48').Value = '  -4.70%  '
$ws.Range('E50').Value = '  -4.06%  '
$ws.Range('E51').Value = '  -3.01%  '

$wb.Save()
